$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "1m89"
$ws.Range("E23").Value = "1m81"

$ws.Range("C30").Select()
